$wb = $excel.ActiveWorkbook

# Overview sheet: row 3 corresponds to the bfbf0c67...md file (de-de language column)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-32-17 06:32:08"

# zh-cn sheet: row 3 is the bfbf0c67 entry - status and handoff datetime updated
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-17 06:32:00"

# de-de sheet: row 3 is the bfbf0c67 entry - status and handoff datetime updated
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-17 06:32:08"
